$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new text cells in the same order the original author
# typed them (column by column, within each 1-3 row block), so the
# rebuilt shared-string table lines up with the target index order. ---
$ws.Range('A5').Value = 'SCRIPT/T01P01A/um0601.ssb'
$ws.Range('C5').Value = ' I\''ve heard that the Time Gear of\n[CS:P]Treeshroud Forest[CR] was stolen.'
$ws.Range('C6').Value = ' Because of that, time has\nstopped in [CS:P]Treeshroud Forest[CR]…'
$ws.Range('C7').Value = ' There\''s been so much horrible\nnews lately. It\''s discouraging.'
$ws.Range('D5').Value = ' Я слышал, что кто-то украл\nШестерню Времени [CS:P]Заросшего Леса[CR].'
$ws.Range('D6').Value = ' Из-за этого, в [CS:P]Заросшем Лесу[CR]\nостановилось время...'
$ws.Range('D7').Value = ' В последнее время происходит\nстолько всего плохого. Это расстраивает.'
$ws.Range('E5').Value = ' Ÿ òìúšàì, œóï ëóï-óï ôëñàì\nŠåòóåñîý Âñåíåîé [CS:P]Èàñïòšåãï Ìåòà[CR].'
$ws.Range('E6').Value = ' Éè-èà üóïãï, â [CS:P]Èàñïòšåí Ìåòô[CR]\nïòóàîïâéìïòû âñåíÿ...'
$ws.Range('E7').Value = ' Â ðïòìåäîåå âñåíÿ ðñïéòöïäéó\nòóïìûëï âòåãï ðìïöïãï. Üóï ñàòòóñàéâàåó.'
$ws.Range('A8').Value = 'SCRIPT/T01P01A/um0605.ssb'
$ws.Range('C8').Value = ' There\''s been a terrible stench\nlately. Or am I imagining it?'
$ws.Range('D8').Value = ' В последнее время я чувствую\nужасную вонь. Или мне кажется?'
$ws.Range('E8').Value = ' Â ðïòìåäîåå âñåíÿ ÿ œôâòóâôý\nôçàòîôý âïîû. Éìé íîå ëàçåóòÿ?'
$ws.Range('C9').Value = ' I\''ve heard.[K] [CS:N]Wigglytuff[CR]\''s Guild is\ngoing on an expedition soon?'
$ws.Range('A9').Value = 'SCRIPT/T01P01A/um0705.ssb'
$ws.Range('C10').Value = ' Good luck with that!'
$ws.Range('C11').Value = ' I hope you get picked for\nthe expedition!'
$ws.Range('D9').Value = ' Я всё слышал.[K] Гильдия\n[CS:N]Виглитаффа[CR] снаряжает экспедицию?'
$ws.Range('D10').Value = ' Я желаю вам удачи!'
$ws.Range('D11').Value = ' Надеюсь, вас возьмут в\nэкспедицию!'
$ws.Range('E9').Value = ' Ÿ âòæ òìúšàì.[K] Ãéìûäéÿ\n[CS:N]Âéãìéóàõõà[CR] òîàñÿçàåó üëòðåäéøéý?'
$ws.Range('E10').Value = ' Ÿ çåìàý âàí ôäàœé!'
$ws.Range('E11').Value = ' Îàäåýòû, âàò âïèûíôó â\nüëòðåäéøéý!'
$ws.Range('C12').Value = ' I\''ve heard! You were picked for\nthe expedition!'
$ws.Range('C13').Value = ' Congratulations!'
$ws.Range('A12').Value = 'SCRIPT/T01P01A/um0801.ssb'
$ws.Range('D12').Value = ' Я всё знаю! Вас взяли в\nэкспедицию!'
$ws.Range('D13').Value = ' Поздравляю!'
$ws.Range('E12').Value = ' Ÿ âòæ èîàý! Âàò âèÿìé â\nüëòðåäéøéý!'
$ws.Range('E13').Value = ' Ðïèäñàâìÿý!'

# --- Fill in the numeric "line number" column (column B) for each new row. ---
$ws.Range('B5').Value = 289
$ws.Range('B6').Value = 292
$ws.Range('B7').Value = 295
$ws.Range('B8').Value = 270
$ws.Range('B9').Value = 245
$ws.Range('B10').Value = 248
$ws.Range('B11').Value = 251
$ws.Range('B12').Value = 223
$ws.Range('B13').Value = 226

# --- Add the thin separator borders that close off each block of rows. ---
# Row 4 gets a bottom border (it now closes off the first, pre-existing block).
$ws.Range('A4:E4').Borders.Item(9).LineStyle = 1
# Rows 5-7 are one block; put a bottom border under row 7.
$ws.Range('A7:E7').Borders.Item(9).LineStyle = 1
# Row 8 is a lone-row block; give it both a top and a bottom border.
$ws.Range('A8:E8').Borders.Item(8).LineStyle = 1
$ws.Range('A8:E8').Borders.Item(9).LineStyle = 1
# Rows 9-11 are one block; put a bottom border under row 11.
$ws.Range('A11:E11').Borders.Item(9).LineStyle = 1
# Rows 12-13 are one block; put a bottom border under row 13.
$ws.Range('A13:E13').Borders.Item(9).LineStyle = 1

# --- Row heights (matches the auto-fit heights Excel computed for the wrapped text). ---
$ws.Rows.Item(4).RowHeight = 21.6
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 21.6
$ws.Rows.Item(7).RowHeight = 21.6
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 21.6
$ws.Rows.Item(12).RowHeight = 43.2

# --- Final selection / scroll position, matching the author leaving off mid-edit. ---
$ws.Range('C8').Select()
